$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (CasesTab): query text changed (now the "Study Code / Age" query),
# StatQuery text unchanged but now shared with rows 3-5; row height grows
# from 270 -> 285.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['West Highland White Terrier'] `nMATCH (c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (co:cohort)<-[*]-(c)`nWITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age`nRETURN  `n       coalesce(c.case_id, '') AS ``Case ID``,`n       coalesce(s.clinical_study_designation, '') AS ``Study Code``,`n       coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n       coalesce(demo.breed, '') AS Breed ,`n       coalesce(diag.disease_term, '') AS Diagnosis ,`n       coalesce(diag.stage_of_disease, '') AS ``Stage of Disease``,`n       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,`n       coalesce(demo.sex, '') AS Sex,`n       coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n       coalesce(demo.weight, '') AS ``Weight (kg)``,`n       coalesce(diag.best_response, '') AS ``Response to Treatment``,`n       coalesce(co.cohort_description, '') AS ``Cohort``"
$ws.Range("C2").Value = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['West Highland White Terrier'] `nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Rows.Item(2).RowHeight = 285

# ---------------------------------------------------------------------------
# Row 3 (SamplesTab): query text lightly reworded (trailing-space cleanup
# around Breed/Diagnosis columns); row height unchanged (225).
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) `nWHERE demo.breed IN ['West Highland White Terrier'] `nWITH DISTINCT samp AS samp, c, demo, diag`nRETURN  coalesce(samp.sample_id, '') AS ``Sample ID``, `n        coalesce(c.case_id, '') AS ``Case ID``, `n        coalesce(demo.breed,'') AS Breed,`n        coalesce(diag.disease_term,'') AS Diagnosis, `n        coalesce(samp.sample_site, '') AS ``Sample Site``,`n        coalesce(samp.summarized_sample_type, '') AS ``Sample Type``,`n        coalesce(samp.specific_sample_pathology, '') AS ``Pathology/Morphology``,`n        coalesce(samp.tumor_grade, '') AS ``Tumor Grade``,`n        coalesce(samp.sample_chronology, '') AS ``Sample Chronology``,`n        coalesce(samp.percentage_tumor, '') AS ``Percentage Tumor``,`n        coalesce(samp.necropsy_sample, '') AS ``Necropsy Sample``,`n        coalesce(samp.sample_preservation, '') AS ``Sample Preservation``"
$ws.Range("C3").Value = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['West Highland White Terrier'] `nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"

# ---------------------------------------------------------------------------
# Row 4 (FilesTab): query gains File Type / byte-size formatting / sample
# join; D4/E4 now point at the TC46 (not TC40) file names; row grows to the
# maximum Excel row height (409.5).
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f)-[*]->(samp:sample)`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['West Highland White Terrier'] `nOPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)`nWITH`n        f, parent, c, demo, diag, s, samp,`n        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n        toInteger(floor(log(f.file_size)/log(1024))) as i,`n        2 as precision`nWITH`n        f, parent, c, demo, diag, s, samp,`n        f.file_size /(1024^i) AS value, `n        10^precision AS factor,`n        units[i] as unit`nWITH    `n        f, parent, c, demo, diag, s, samp, unit,`n        round(factor * value)/factor AS size`nRETURN `n        coalesce(f.file_name, '') AS ``File Name``,`n        coalesce(f.file_type, '') AS ``File Type``,`n        coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n        coalesce(samp.sample_id, '') AS ``Sample ID``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n        coalesce(demo.breed,'') AS Breed ,`n        coalesce(diag.disease_term,'') AS Diagnosis"
$ws.Range("C4").Value = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['West Highland White Terrier'] `nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("D4").Value = "TC46_Canine_Filter_Breed-WestHlnd_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC46_Canine_Filter_Breed-WestHlnd_WebData.xlsx"
$ws.Rows.Item(4).RowHeight = 409.5

# ---------------------------------------------------------------------------
# Row 5 (new StudyFilesTab): entirely new row.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = "MATCH (f:file)-->(s:study)<--(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['West Highland White Terrier']`nWITH`n    f, s,`n    ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,`n    toInteger(floor(log(f.file_size)/log(1024))) as i,`n    2 as precision`nWITH`n    f, s,`n    f.file_size /(1024^i) AS value, 10^precision AS factor,`n    units[i] as unit`nWITH`n    f, s, unit,`n    round(factor * value)/factor AS size`nRETURN DISTINCT`n  coalesce(f.file_name, '') AS ``File Name``,`n  coalesce(f.file_type, '') AS ``File Type``,`n  coalesce(`"study`", '') AS ``Association``,`n  coalesce(f.file_description, '') AS ``Description``,`n  coalesce(f.file_format, '') AS  Format,`n  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,`n  coalesce(s.clinical_study_designation,'') AS ``Study Code``"
$ws.Range("C5").Value = "MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)`nOPTIONAL MATCH (samp:sample)-->(c)`nOPTIONAL MATCH (diag:diagnosis)-->(c)`nOPTIONAL MATCH (f:file)-[*]->(c)`nOPTIONAL MATCH (sf:file)-->(s)`nWITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p`nWHERE demo.breed IN ['West Highland White Terrier'] `nRETURN  `n    count(distinct p) AS Programs,`n    count(distinct s) AS Studies,`n    count(distinct c) AS Cases,`n    count(distinct samp) AS Samples,`n    count(distinct f) AS ``Case Files``,`n    count(distinct sf) AS ``Study Files``"
$ws.Range("D5").Value = "TC46_Canine_Filter_Breed-WestHlnd_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC46_Canine_Filter_Breed-WestHlnd_WebData.xlsx"

# Same wrap-text cell style used by B2:C4.
$ws.Range("B5:C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 330

# ---------------------------------------------------------------------------
# View state: window scrolled/zoomed to show the newly added row, and the
# active selection moved onto it.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 55
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B5").Select()
